$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("openbis-metadata")

# Insert a new row above the current "Timepoint Type" row (row 4) to hold
# the new "Reference Strain" property.
$ws.Rows.Item(4).Insert()

# Fill in the new "Reference Strain" row (B4 stays empty).
$ws.Range("A4").Value = "Reference Strain"
$ws.Range("C4").Value = "The Reference Strain (for relative quantification data sets, leave empty for absolute)"

# The inserted row picked up plain formatting for C4; restore the left
# border used by the rest of the Description column cells (style index 2).
$ws.Range("C4").Borders.Item(7).LineStyle = 1

# Update the "Value Unit" row (now row 7): simplify the unit text and its
# description (drop the "protein digest" / stray "or" wording).
$ws.Range("B7").Value = "fmol/ug"
$ws.Range("C7").Value = "One of mM, uM, Percent, RatioT1, RatioCs, AU, Dimensionless, fmol/ug"

# Clear the "Scale" row's description (now row 8) - it becomes empty.
$ws.Range("C8").Value = ""

# Widen column A slightly to fit the new "Reference Strain" label.
$ws.Columns.Item(1).ColumnWidth = 17.57

# Restore the (somewhat arbitrary) saved selection state from the source file.
$ws.Range("C16").Select() | Out-Null
